# Scheduled-runner update of market/profit figures (columns H-N:
# currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR
# sheets. Each block below is keyed by the row's "Leve Item ID" (column G)
# for traceability. Cells that the source feed no longer produces a number
# for are cleared outright (matching rows losing their <c> element);
# cells that newly gained a number are written for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 152.81818
$ws.Range("I4").Value = 193
$ws.Range("K4").Value = 193
$ws.Range("M4").Value = -79

# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 1897.25
$ws.Range("I12").Value = 1899.625
$ws.Range("J12").Value = 1892.5
$ws.Range("K12").Value = 1899.625
$ws.Range("L12").Value = 1892.5
$ws.Range("M12").Value = -1729.625
$ws.Range("N12").Value = -2232.5

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 5833.3335
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 778.25
$ws.Range("I28").Value = 662.35297
$ws.Range("J28").Value = 1435
$ws.Range("K28").Value = 662.35297
$ws.Range("L28").Value = 1435
$ws.Range("M28").Value = -177.35297
$ws.Range("N28").Value = -2405

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 11115267
$ws.Range("I40").Value = 33335666
$ws.Range("J40").Value = 5066.6665
$ws.Range("K40").Value = 33335666
$ws.Range("L40").Value = 5066.6665
$ws.Range("M40").Value = -33335491
$ws.Range("N40").Value = -5416.6665

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 6138.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6138.5
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 6138.5
$ws.Range("N86").Value = -8384.5

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 6138.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6138.5
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 30692.5
$ws.Range("N89").Value = -41924.5

# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 2589.3845
$ws.Range("J112").Value = 3710.4285
$ws.Range("L112").Value = 11131.2855
$ws.Range("N112").Value = -13347.2855

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 3344
$ws.Range("I141").Value = 1602.5
$ws.Range("K141").Value = 4807.5
$ws.Range("M141").Value = 372.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1486.6666
$ws.Range("I2").Value = 1057.4615
$ws.Range("K2").Value = 1057.4615
$ws.Range("M2").Value = -944.4614999999999

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 4864.769
$ws.Range("I45").Value = 3474.125
$ws.Range("K45").Value = 3474.125
$ws.Range("M45").Value = -3097.125

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1486.6666
$ws.Range("I116").Value = 1057.4615
$ws.Range("K116").Value = 1057.4615
$ws.Range("M116").Value = 1236.5385

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 997.5
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 995
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 2985
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -7885

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1486.6666
$ws.Range("I3").Value = 1057.4615
$ws.Range("K3").Value = 1057.4615
$ws.Range("M3").Value = -943.4614999999999

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1547
$ws.Range("I86").Value = 1404.5294
$ws.Range("J86").Value = 2152.5
$ws.Range("K86").Value = 1404.5294
$ws.Range("L86").Value = 2152.5
$ws.Range("M86").Value = -281.5293999999999
$ws.Range("N86").Value = -4398.5

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1547
$ws.Range("I89").Value = 1404.5294
$ws.Range("J89").Value = 2152.5
$ws.Range("K89").Value = 7022.646999999999
$ws.Range("L89").Value = 10762.5
$ws.Range("M89").Value = -1406.646999999999
$ws.Range("N89").Value = -21994.5

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 2306.4375
$ws.Range("I99").Value = 2160.0833
$ws.Range("K99").Value = 2160.0833
$ws.Range("M99").Value = -662.0832999999998

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2976.3157
$ws.Range("I105").Value = 2598.7273
$ws.Range("K105").Value = 2598.7273
$ws.Range("M105").Value = -851.7273

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 10149.75
$ws.Range("I107").Value = 11866.333
$ws.Range("K107").Value = 11866.333
$ws.Range("M107").Value = -9946.333000000001

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2307.7615
$ws.Range("I134").Value = 1324.4595
$ws.Range("J134").Value = 3021.1372
$ws.Range("K134").Value = 3973.3785
$ws.Range("L134").Value = 9063.411599999999
$ws.Range("M134").Value = -1438.3785
$ws.Range("N134").Value = -14133.4116

$ws = $wb.Worksheets.Item("CRP")
# Row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 4218.5454
$ws.Range("I41").Value = 4218.5454
$ws.Range("K41").Value = 4218.5454
$ws.Range("M41").Value = -3790.5454

# Row 54 (Leve Item ID 2413)
$ws.Range("H54").Value = 50192
$ws.Range("J54").Value = 50192
$ws.Range("L54").Value = 50192
$ws.Range("N54").Value = -51508

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1567
$ws.Range("I107").Value = 655.3333
$ws.Range("J107").Value = 2934.5
$ws.Range("K107").Value = 655.3333
$ws.Range("L107").Value = 2934.5
$ws.Range("M107").Value = 1264.6667
$ws.Range("N107").Value = -6774.5

# Row 116 (Leve Item ID 26117)
$ws.Range("H116").Value = 89000
$ws.Range("J116").Value = 89000
$ws.Range("L116").Value = 89000
$ws.Range("N116").Value = -98178

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1758.1666
$ws.Range("I134").Value = 1626.7435
$ws.Range("K134").Value = 4880.2305
$ws.Range("M134").Value = -2345.2305

$ws = $wb.Worksheets.Item("CUL")
# Row 44 (Leve Item ID 4702)
$ws.Range("H44").Value = 679.4
$ws.Range("I44").Value = 584.8570999999999
$ws.Range("J44").Value = 900
$ws.Range("K44").Value = 1754.5713
$ws.Range("L44").Value = 2700
$ws.Range("M44").Value = -1356.5713
$ws.Range("N44").Value = -3496

# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 725.0909
$ws.Range("I117").Value = 666.3333
$ws.Range("J117").Value = 747.125
$ws.Range("K117").Value = 1998.9999
$ws.Range("L117").Value = 2241.375
$ws.Range("M117").Value = 1443.0001
$ws.Range("N117").Value = -9125.375

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1152.96
$ws.Range("I131").Value = 659.0769
$ws.Range("K131").Value = 1977.2307
$ws.Range("M131").Value = 3062.7693

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4319.5
$ws.Range("I80").Value = 4573.625
$ws.Range("K80").Value = 4573.625
$ws.Range("M80").Value = -3575.625

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4319.5
$ws.Range("I83").Value = 4573.625
$ws.Range("K83").Value = 22868.125
$ws.Range("M83").Value = -17876.125

# Row 118 (Leve Item ID 26172)
$ws.Range("H118").Value = 81061
$ws.Range("I118").Value = 50000
$ws.Range("J118").Value = 88826.25
$ws.Range("K118").Value = 50000
$ws.Range("L118").Value = 88826.25
$ws.Range("M118").Value = -48343
$ws.Range("N118").Value = -92140.25

# Row 141 (Leve Item ID 42504)
$ws.Range("H141").Value = 139714.5
$ws.Range("J141").Value = 139714.5
$ws.Range("L141").Value = 139714.5
$ws.Range("N141").Value = -150074.5

$ws = $wb.Worksheets.Item("LTW")
# Row 5 (Leve Item ID 3790)
$ws.Range("H5").Value = 11839.333
$ws.Range("I5").Value = 5259
$ws.Range("K5").Value = 5259
$ws.Range("M5").Value = -5146

# Row 20 (Leve Item ID 4308)
$ws.Range("H20").Value = 10324.667
$ws.Range("J20").Value = 19969
$ws.Range("L20").Value = 19969
$ws.Range("N20").Value = -20421

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 8226.478999999999
$ws.Range("I40").Value = 7496.6313
$ws.Range("J40").Value = 11693.25
$ws.Range("K40").Value = 7496.6313
$ws.Range("L40").Value = 11693.25
$ws.Range("M40").Value = -7360.6313
$ws.Range("N40").Value = -11965.25

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 13308
$ws.Range("I46").Value = 5957.3
$ws.Range("J46").Value = 34310
$ws.Range("K46").Value = 5957.3
$ws.Range("L46").Value = 34310
$ws.Range("M46").Value = -5769.3
$ws.Range("N46").Value = -34686

$ws = $wb.Worksheets.Item("WVR")
# Row 21 (Leve Item ID 3341)
$ws.Range("H21").Value = 19507.5
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 23 (Leve Item ID 3325)
$ws.Range("H23").Value = 1309
$ws.Range("J23").Value = 900
$ws.Range("L23").Value = 900
$ws.Range("N23").Value = -1358

# Row 24 (Leve Item ID 3561)
$ws.Range("H24").Value = 19504.5
$ws.Range("I24").Value = 19504.5
$ws.Range("K24").Value = 19504.5
$ws.Range("M24").Value = -19274.5

# Row 25 (Leve Item ID 3064)
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15586

# Row 30 (Leve Item ID 2700)
$ws.Range("H30").Value = 5504.5
$ws.Range("I30").Value = 6009
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 6009
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -5902
$ws.Range("N30").Value = -5214

# Row 35 (Leve Item ID 3341)
$ws.Range("H35").Value = 19507.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 74 (Leve Item ID 19022)
$ws.Range("H74").Value = 16093.5
$ws.Range("I74").Value = 20523
$ws.Range("J74").Value = 11664
$ws.Range("K74").Value = 20523
$ws.Range("L74").Value = 11664
$ws.Range("M74").Value = -19587
$ws.Range("N74").Value = -13536

# Row 77 (Leve Item ID 19022)
$ws.Range("H77").Value = 16093.5
$ws.Range("I77").Value = 20523
$ws.Range("J77").Value = 11664
$ws.Range("K77").Value = 61569
$ws.Range("L77").Value = 34992
$ws.Range("M77").Value = -56889
$ws.Range("N77").Value = -44352

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 46319.707
$ws.Range("I81").Value = 69918.60000000001
$ws.Range("J81").Value = 6988.222
$ws.Range("K81").Value = 139837.2
$ws.Range("L81").Value = 13976.444
$ws.Range("M81").Value = -138776.2
$ws.Range("N81").Value = -16098.444

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 46319.707
$ws.Range("I84").Value = 69918.60000000001
$ws.Range("J84").Value = 6988.222
$ws.Range("K84").Value = 699186
$ws.Range("L84").Value = 69882.22
$ws.Range("M84").Value = -693882
$ws.Range("N84").Value = -80490.22

# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1805.3
$ws.Range("J100").Value = 1867.6666
$ws.Range("L100").Value = 3735.3332
$ws.Range("N100").Value = -4817.3332

# Row 120 (Leve Item ID 26310)
$ws.Range("H120").Value = 91273.336
$ws.Range("J120").Value = 91273.336
$ws.Range("L120").Value = 91273.336
$ws.Range("N120").Value = -100949.336
